$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''41.650.57'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +4.86%  '
$ws.Cells.Item(3, 4).Value = '''2.228.65'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +3.27%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '''228.82'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.73%  '
$ws.Cells.Item(6, 4).Value = '''0.623'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.63%  '
$ws.Cells.Item(7, 4).Value = '''61.46'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -2.76%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).Value = '''0.402'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.74%  '
$ws.Cells.Item(10, 4).Value = '''58.04'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.52%  '
$ws.Cells.Item(11, 4).Value = '''0.0879'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +3.94%  '
$ws.Cells.Item(12, 5).Value = '  +0.04%  '
$ws.Cells.Item(13, 4).Value = '''2.555.28'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +3.01%  '
$ws.Cells.Item(14, 4).Value = '''15.67'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -1.34%  '
$ws.Cells.Item(15, 4).Value = '''21.57'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -1.18%  '
$ws.Cells.Item(16, 4).Value = '''0.795'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.92%  '
$ws.Cells.Item(17, 4).Value = '''5.56'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.56%  '
$ws.Cells.Item(18, 4).Value = '''2.226.58'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +3.25%  '
$ws.Cells.Item(19, 4).Value = '''41.624.31'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +5.08%  '
$ws.Cells.Item(20, 4).Value = '''72.67'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.45%  '
$ws.Cells.Item(21, 4).Value = '''0.0₃0889'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +5.52%  '
$ws.Cells.Item(22, 4).Value = '''6.04'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.48%  '
$ws.Cells.Item(23, 4).Value = '''248.05'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +7.84%  '
$ws.Cells.Item(24, 5).Value = '  +0.01%  '
$ws.Cells.Item(25, 4).Value = '''2.37'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.12%  '
$ws.Cells.Item(26, 4).Value = '''2.27'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.66%  '
$ws.Cells.Item(27, 4).Value = '''9.48'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.58%  '
$ws.Cells.Item(28, 4).Value = '''167.87'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -2.54%  '
$ws.Cells.Item(29, 5).Value = '  +0.40%  '
$ws.Cells.Item(30, 4).Value = '''19.94'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.51%  '
$ws.Cells.Item(31, 4).Value = '''1.42'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.90%  '
$ws.Cells.Item(32, 4).Value = '''2.60'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.36%  '
$ws.Cells.Item(33, 5).Value = '  -0.04%  '
$ws.Cells.Item(34, 5).Value = '  +7.90%  '
$ws.Cells.Item(35, 4).Value = '''4.65'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.60%  '
$ws.Cells.Item(36, 4).Value = '''0.0625'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.16%  '
$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).Value = '''3.70'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.37%  '
$ws.Cells.Item(38, 2).Value = 'THORChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(38, 4).Value = '''6.57'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -5.31%  '
$ws.Cells.Item(39, 4).Value = '''2.38'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.58%  '
$ws.Cells.Item(40, 4).Value = '''0.999'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -0.08%  '
$ws.Cells.Item(41, 4).Value = '''0.000237'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +28.99%  '
$ws.Cells.Item(42, 4).Value = '''4.86'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -3.92%  '
$ws.Cells.Item(43, 4).Value = '''0.0236'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +4.77%  '
$ws.Cells.Item(44, 4).Value = '''8.74'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +12.93%  '
$ws.Cells.Item(45, 4).Value = '''0.0984'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +7.03%  '
$ws.Cells.Item(46, 4).Value = '''99.08'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.50%  '
$ws.Cells.Item(47, 4).Value = '''1.470.55'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.95%  '
$ws.Cells.Item(48, 5).Value = '  -2.41%  '
$ws.Cells.Item(49, 4).Value = '''16.47'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -6.58%  '
$ws.Cells.Item(50, 5).Value = '  -0.85%  '
$ws.Cells.Item(51, 5).Value = '  -0.76%  '
